# edit.ps1 - applies the "Top Egresos" sheet addition + Resumen cleanup
# described by the commit: adds a new "Top Egresos" worksheet (ranking of the
# 10 largest egresos) between "Egresos" and "Prestadores", and removes the
# now-redundant "Pacientes Transferencia" (row 24) / "Sueldos" (row 32) rows
# from the "Resumen" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Resumen sheet: drop "Pacientes Transferencia" and "Sueldos" rows.
#    Deleting row 24 shifts everything below it up by one, so the old row 32
#    ("Sueldos") becomes row 31; delete that too so the sheet ends at B30.
# ---------------------------------------------------------------------------
$resumen = $wb.Worksheets.Item("Resumen")
$resumen.Rows.Item(24).Delete()
$resumen.Rows.Item(31).Delete()

# ---------------------------------------------------------------------------
# 2) Insert the new "Top Egresos" sheet right after "Egresos" (i.e. before
#    "Prestadores").
# ---------------------------------------------------------------------------
$egresos = $wb.Worksheets.Item("Egresos")
$topEgresos = $wb.Worksheets.Add($null, $egresos)
$topEgresos.Name = "Top Egresos"

# Column widths (chars): A=9 B=21 C=23 D=50 E=10 F=14 G=13
# (ColumnWidth is expressed in character units; this engine's Width<->chars
# conversion carries a constant +5/6 offset, so subtract it back out.)
$widths = @(9, 21, 23, 50, 10, 14, 13)
for ($i = 0; $i -lt $widths.Length; $i++) {
    $topEgresos.Columns.Item($i + 1).ColumnWidth = $widths[$i] - 0.8333333333
}

# Header row
$headers = @("Ranking", "Fecha", "Concepto", "Detalle", "Monto", "Subcategoría", "Banco")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $topEgresos.Cells.Item(1, $i + 1).Value = $headers[$i]
}
$hdrRange = $topEgresos.Range("A1:G1")
$hdrRange.Font.Bold = $true
$hdrRange.Font.Color = 16777215
$hdrRange.Interior.Color = 12874308
$hdrRange.HorizontalAlignment = -4108
$hdrRange.VerticalAlignment = -4108
$hdrRange.Borders.LineStyle = 1

# Data rows: Ranking, Fecha (serial datetime), Concepto, Detalle, Monto, Subcategoria, Banco
$rows = @(
    @(1, 45945.77690972222, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 1500000, "Prestadores", "Supervielle"),
    @(2, 45939.78494212963, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 862300, "Prestadores", "Supervielle"),
    @(3, 45954.57414351852, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 600000, "Prestadores", "Supervielle"),
    @(4, 45954.57196759259, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 600000, "Prestadores", "Supervielle"),
    @(5, 45954.57053240741, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 600000, "Prestadores", "Supervielle"),
    @(6, 45950.68238425926, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 375000, "Prestadores", "Supervielle"),
    @(7, 45939.82627314814, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 350000, "Prestadores", "Supervielle"),
    @(8, 45950.68052083333, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 304500, "Prestadores", "Supervielle"),
    @(9, 45938.59410879629, "Compra Visa Débito",     "COMERCIO: PAGOS360*EPEC OPERACION: 139928",             239243.1, "Servicios",  "Supervielle"),
    @(10, 45938.48383101852, "Transferencia por CBU", "CONCEPTO: Transferencia enviada TERMINAL: TESP0000...", 204254, "Prestadores", "Supervielle")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowNum = $r + 2
    $data = $rows[$r]
    $topEgresos.Cells.Item($rowNum, 1).Value = $data[0]
    $topEgresos.Cells.Item($rowNum, 2).Value = $data[1]
    $topEgresos.Cells.Item($rowNum, 3).Value = $data[2]
    $topEgresos.Cells.Item($rowNum, 4).Value = $data[3]
    $topEgresos.Cells.Item($rowNum, 5).Value = $data[4]
    $topEgresos.Cells.Item($rowNum, 6).Value = $data[5]
    $topEgresos.Cells.Item($rowNum, 7).Value = $data[6]
}

# Date/time format for the "Fecha" column (matches the rest of the workbook)
$topEgresos.Range("B2:B11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

Write-Host "Top Egresos sheet created and Resumen rows cleaned up"
